$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 12:52"

# 2) Update country data rows with refreshed figures.
#    Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#             E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 5: Espana
$ws.Range("B5").Value = 219764
$ws.Range("C5").Value = 6740
$ws.Range("D5").Value = 92355
$ws.Range("E5").Value = 104885
$ws.Range("G5").Value = 367
$ws.Range("H5").Value = 22524

# Row 8: Alemania
$ws.Range("B8").Value = 153307
$ws.Range("C8").Value = 178
$ws.Range("E8").Value = 40932

# Rows 37-40: Catar moves above Indonesia/Dinamarca/Bielorrusia in the
# country ordering (alphabetically-independent ranking list), and Catar's
# figures get refreshed while Indonesia/Dinamarca/Bielorrusia keep their
# previous figures but shift down one row.
$ws.Range("A37").Value = "Catar"
$ws.Range("B37").Value = 8525
$ws.Range("C37").Value = 761
$ws.Range("D37").Value = 750
$ws.Range("E37").Value = 7765
$ws.Range("F37").Value = 72
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 10

$ws.Range("A38").Value = "Indonesia"
$ws.Range("B38").Value = 8211
$ws.Range("C38").Value = 436
$ws.Range("D38").Value = 1002
$ws.Range("E38").Value = 6520
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 42
$ws.Range("H38").Value = 689

$ws.Range("A39").Value = "Dinamarca"
$ws.Range("B39").Value = 8210
$ws.Range("C39").Value = 137
$ws.Range("D39").Value = 5384
$ws.Range("E39").Value = 2432
$ws.Range("F39").Value = 74
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 394

$ws.Range("A40").Value = "Bielorrusia"
$ws.Range("B40").Value = 8022
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 938
$ws.Range("E40").Value = 7024
$ws.Range("F40").Value = 92
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 60

# Row 76: Bosnia y Herzegovina
$ws.Range("B76").Value = 1421
$ws.Range("C76").Value = 8
$ws.Range("D76").Value = 538
$ws.Range("E76").Value = 828
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 55

# Row 107: Malta
$ws.Range("B107").Value = 447
$ws.Range("C107").Value = 2
$ws.Range("D107").Value = 223
$ws.Range("E107").Value = 221

# Row 190: San Cristobal y Nieves
$ws.Range("D190").Value = 2
$ws.Range("E190").Value = 13
